$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Aperture" column (F) values for rows 2-16: from 1E-4 to 1,
# and apply the scientific-number style (same style/format used by column E).
$range = $ws.Range("F2:F16")
$range.Value = 1
$range.NumberFormat = "0.00E+00"

# Update the selection on the sheet to Q10 (single cell).
$ws.Range("Q10").Select()
